$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Revisions to existing rows (column F = AgTests, column G = AgPosit) ---
$updates = @(
    @{ Row = 314; F = 65186; G = 3134 },
    @{ Row = 320; F = 73998; G = 3358 },
    @{ Row = 322; F = 109725 },
    @{ Row = 324; F = 249422 },
    @{ Row = 325; F = 769893 },
    @{ Row = 334; F = 193279 },
    @{ Row = 336; F = 82036; G = 2573 },
    @{ Row = 341; F = 283701 },
    @{ Row = 349; F = 159665; G = 2756 },
    @{ Row = 435; F = 83137 },
    @{ Row = 436; F = 139211 },
    @{ Row = 439; F = 86643 },
    @{ Row = 440; F = 72863 },
    @{ Row = 441; F = 65807 },
    @{ Row = 442; F = 67286 },
    @{ Row = 443; F = 102776 },
    @{ Row = 449; F = 59728 },
    @{ Row = 455; F = 49993 },
    @{ Row = 456; F = 47796 },
    @{ Row = 464; F = 69788 },
    @{ Row = 469; F = 38996 },
    @{ Row = 473; F = 38136 },
    @{ Row = 474; F = 43309 },
    @{ Row = 475; F = 33242 },
    @{ Row = 476; F = 34186; G = 29 },
    @{ Row = 477; F = 35304; G = 32 }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 6).Value = $u.F
    if ($u.ContainsKey("G")) {
        $ws.Cells.Item($u.Row, 7).Value = $u.G
    }
}

# --- Append new daily rows 478-480 ---
$newRows = @(
    @{ Row = 478; A = 44372; B = 391531; C = 5685; D = 41; E = 12505; F = 44687; G = 26 },
    @{ Row = 479; A = 44373; B = 391551; C = 2553; D = 20; E = 12505; F = 33006; G = 28 },
    @{ Row = 480; A = 44374; B = 391566; C = 804;  D = 15; E = 12509; F = 25563; G = 18 }
)

foreach ($r in $newRows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.A
    $ws.Cells.Item($r.Row, 1).NumberFormat = "yyyy-mm-dd"
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $ws.Cells.Item($r.Row, 4).Value = $r.D
    $ws.Cells.Item($r.Row, 5).Value = $r.E
    $ws.Cells.Item($r.Row, 6).Value = $r.F
    $ws.Cells.Item($r.Row, 7).Value = $r.G
}
